# Fix title of contribution form:
# "Assignment 1: Project Plan" -> "Assignment 2: Machine Learning"
# (the NBSP between "Project" and "Plan" in the source text is matched
#  via wildcard so we don't have to worry about exact whitespace char)

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Assignment 1: Project*Plan",  # FindText
    $false,                         # MatchCase
    $false,                         # MatchWholeWord
    $true,                          # MatchWildcards
    $false,                         # MatchSoundsLike
    $false,                         # MatchAllWordForms
    $true,                          # Forward
    1,                              # Wrap (wdFindContinue)
    $false,                         # Format
    "Assignment 2: Machine Learning", # ReplaceWith
    2                               # Replace (wdReplaceAll)
)
